$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended below the existing data (row 5)
$row = 5

$ws.Cells.Item($row, 1).Value = 42611.887627314813
$ws.Cells.Item($row, 2).Value = 12
$ws.Cells.Item($row, 3).Value = 55
$ws.Cells.Item($row, 4).Value = 42
$ws.Cells.Item($row, 5).Value = 100
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 8269
$ws.Cells.Item($row, 8).Value = 7214
$ws.Cells.Item($row, 9).Value = 1214
$ws.Cells.Item($row, 10).Value = 143
$ws.Cells.Item($row, 11).Value = 109
$ws.Cells.Item($row, 12).Value = 1
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = "Noun"

# Column A keeps the same date/time number format as the rows above it
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"
